# Update "Prefeitura de Itapema_itaitaitia.xlsx"
#
# 1. Duplicate the "Centro da Cidade" sheet, placing the copy right after it,
#    and rename the copy to "Fundos da Cidade" (it keeps the original
#    "Rena" budget line).
# 2. Replace the data on the (now repurposed) "Centro da Cidade" sheet with
#    the new "Pastor e Ovelha" budget line.
# 3. Add a "Centro da Cidade" scenario row to the "cenarios" sheet, and
#    repoint the existing row 2 at "Fundos da Cidade".
# 4. Switch workbook calculation to manual.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "Centro da Cidade" -> "Fundos da Cidade" -----------------
$centro = $wb.Worksheets.Item("Centro da Cidade")
$centro.Copy($null, $centro)
$fundos = $wb.Worksheets.Item("Centro da Cidade (2)")
$fundos.Name = "Fundos da Cidade"

# --- 2. Replace the data on "Centro da Cidade" ------------------------------
$ws = $wb.Worksheets.Item("Centro da Cidade")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "AD03"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "0,63"
$ws.Range("F2").Value = "0,7"
$ws.Range("G2").ClearContents()
$ws.Range("H2").Value = 3072.29
$ws.Range("I2").ClearContents()
$ws.Range("J2").Value = "Pastor e Ovelha, produzido em alumínio e pintura branca"
$ws.Range("K2").Value = "Pastor e ovelha aramados, medindo aproximadamente o Pastor 1,85m de altura x 0,65m de largura x 0,75m de comprimento e a ovelha  0,45m  de altura x 0,22m de largura e 0,75m de comprimento, produzidos em alumínio, com pintura em esmalte sintético automotivo na cor branca e verniz automotivo brilhante com detalhes de conjuntos com LEDs brancos, fio elétrico branco 2 x 0,5mm², na tensão de 220v. "
$ws.Range("L2").Value = "Pastor e ovelha aramados, figuras tridimensionais, produzidas com arames de alumíno, pintura em epóxi na cor branca, aplicação de lâmpadas de LED brancas, com fio branco ou morno (medida do fio: 2x0,5mm²). Medidas das figuras: mínimo: Pastor 1,80 de altura x 0,60m de largura x 0,70m de comprimento e Ovelha  0,59m de altura x 0,17m de largura x 0,70m de comprimento  / máximo: Pastor 1,90m de altura x 0,70m de largura x 0,80m de comprimento e  Ovelha  0,49m de altura x 0,27m de largura x 0,80m de comprimento."
# M2 keeps its existing formula (=C2*H2) and recalculates automatically.

# --- 3. Update "cenarios" sheet --------------------------------------------
$cenarios = $wb.Worksheets.Item("cenarios")
$cenarios.Range("B2").Value = "Fundos da Cidade"
$cenarios.Range("A3").Value = 3
$cenarios.Range("B3").Value = "Centro da Cidade"

# --- 4. Manual calculation mode ---------------------------------------------
$excel.Calculation = -4135

# Keep "Centro da Cidade" as the active/selected tab (as it was originally).
$ws.Activate()
